$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()
$ws.Range("C55").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
Write-Host "done"
